$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "FAIL"
$ws.Range("D5").Value = "FAIL"

[void]$ws.Range("D2:D5").Select()
